# Fix bugs after changing column order
# - Swap columns E and F (value + formatting) for rows 2..40 on the first
#   worksheet (the data had been mistakenly entered one column early).
# - Add empty "wrap text" formatted cells in column H for rows 1, 3, 33, 34
#   (these rows grew an extra blank column after the shift).
# - Update the active selection to H4 (matches the new layout focus).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$staging = $ws.Range("ZZ1")

for ($r = 2; $r -le 40; $r++) {
    $colE = $ws.Range("E$r")
    $colF = $ws.Range("F$r")

    $colE.Copy()
    $staging.PasteSpecial()

    $colF.Copy()
    $colE.PasteSpecial()

    $staging.Copy()
    $colF.PasteSpecial()

    $staging.Clear()
}

# New blank (but formatted) cells introduced in column H.
$ws.Range("H1").WrapText = $true
$ws.Range("H3").WrapText = $true
$ws.Range("H33").WrapText = $true
$ws.Range("H34").WrapText = $true

# Update the current selection to reflect the new column H.
$ws.Range("H4").Select()
